$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells that are being updated, to preserve exact
# formatting (leading/trailing zeros, dot-grouped numbers) as plain text,
# matching the original inlineStr cell type.
$priceRows = @(2,3,5,6,9,13,14,17,18,19,20,21,27,29,33,35,38,39,40,41,42,43,46,47,48,50,51)
foreach ($r in $priceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.140.02"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "2.444.12"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "583.82"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").Value = "142.92"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "2.438.38"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("D14").Value = "26.49"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "62.041.60"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "2.430.39"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "10.77"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").Value = "7.15"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "326.33"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("D27").Value = "596.72"
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "2.562.22"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "1.90"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").Value = "4.89"
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "0.375"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "153.21"
$ws.Range("E39").Value = "  +4.23%  "
$ws.Range("D40").Value = "18.39"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").Value = "5.27"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").Value = "43.16"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("D43").Value = "1.71"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").Value = "141.91"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").Value = "3.64"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "0.0₆0266"
$ws.Range("E48").Value = "  +18.48%  "
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "0.0519"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "19.89"
$ws.Range("E51").Value = "  -1.13%  "
